# Auto-generated Excel COM-interop edit script
# Applies Phantom_Profits price-refresh updates across BSM, CUL, and GSM sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 14199.324
$ws.Range("I86").Value = 13124.292
$ws.Range("K86").Value = 13124.292
$ws.Range("M86").Value = -12001.292
$ws.Range("H89").Value = 14199.324
$ws.Range("I89").Value = 13124.292
$ws.Range("K89").Value = 65621.45999999999
$ws.Range("M89").Value = -60005.45999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 36740.434
$ws.Range("J34").Value = 39327.285
$ws.Range("L34").Value = 117981.855
$ws.Range("N34").Value = -118149.855
$ws.Range("H62").Value = 9729.714
$ws.Range("J62").Value = 10999.667
$ws.Range("L62").Value = 32999.001
$ws.Range("N62").Value = -34371.001
$ws.Range("H63").Value = 4056
$ws.Range("J63").Value = 8000
$ws.Range("L63").Value = 24000
$ws.Range("N63").Value = -25498
$ws.Range("H64").Value = 873.6667
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H65").Value = 9729.714
$ws.Range("J65").Value = 10999.667
$ws.Range("L65").Value = 98997.003
$ws.Range("N65").Value = -105861.003
$ws.Range("H66").Value = 4056
$ws.Range("J66").Value = 8000
$ws.Range("L66").Value = 72000
$ws.Range("N66").Value = -79488
$ws.Range("H67").Value = 873.6667
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H69").Value = 2072.8572
$ws.Range("J69").Value = 2168.6667
$ws.Range("L69").Value = 6506.000100000001
$ws.Range("N69").Value = -8128.000100000001
$ws.Range("H72").Value = 2072.8572
$ws.Range("J72").Value = 2168.6667
$ws.Range("L72").Value = 19518.0003
$ws.Range("N72").Value = -27630.0003
$ws.Range("H95").Value = 100000
$ws.Range("J95").Value = 100000
$ws.Range("L95").Value = 300000
$ws.Range("N95").Value = -304118
$ws.Range("H97").Value = 1765.3334
$ws.Range("I97").Value = 198.5
$ws.Range("J97").Value = 4899
$ws.Range("K97").Value = 595.5
$ws.Range("L97").Value = 14697
$ws.Range("M97").Value = -99.5
$ws.Range("N97").Value = -15689
$ws.Range("H98").Value = 256
$ws.Range("J98").Value = 312
$ws.Range("L98").Value = 936
$ws.Range("N98").Value = -3932
$ws.Range("H99").Value = 3100.625
$ws.Range("I99").Value = 4140.6665
$ws.Range("J99").Value = 2860.6155
$ws.Range("K99").Value = 12421.9995
$ws.Range("L99").Value = 8581.8465
$ws.Range("M99").Value = -10175.9995
$ws.Range("N99").Value = -13073.8465
$ws.Range("H100").Value = 5777.2
$ws.Range("I100").Value = 299.33334
$ws.Range("J100").Value = 13994
$ws.Range("K100").Value = 898.0000200000001
$ws.Range("L100").Value = 41982
$ws.Range("M100").Value = -87.00002000000006
$ws.Range("N100").Value = -43604
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = $null
$ws.Range("H103").Value = 4990
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").Value = $null
$ws.Range("H106").Value = 18200
$ws.Range("H107").Value = 456.13794
$ws.Range("J107").Value = 461.1111
$ws.Range("L107").Value = 1383.3333
$ws.Range("N107").Value = -5223.3333
$ws.Range("H108").Value = 3249
$ws.Range("I108").Value = 3249
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 9747
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -6867
$ws.Range("N108").Value = $null
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").Value = $null
$ws.Range("H111").Value = 499
$ws.Range("I111").Value = 499
$ws.Range("K111").Value = 1497
$ws.Range("M111").Value = 1570
$ws.Range("H112").Value = 9900.417
$ws.Range("I112").Value = 2413.5
$ws.Range("J112").Value = 11397.8
$ws.Range("K112").Value = 7240.5
$ws.Range("L112").Value = 34193.39999999999
$ws.Range("M112").Value = -6132.5
$ws.Range("N112").Value = -36409.39999999999
$ws.Range("H114").Value = 1992
$ws.Range("I114").Value = 1215.6666
$ws.Range("J114").Value = 2457.8
$ws.Range("K114").Value = 3646.9998
$ws.Range("L114").Value = 7373.400000000001
$ws.Range("M114").Value = -392.9998000000001
$ws.Range("N114").Value = -13881.4
$ws.Range("H115").Value = 2114
$ws.Range("I115").Value = 921
$ws.Range("J115").Value = 4500
$ws.Range("K115").Value = 2763
$ws.Range("L115").Value = 13500
$ws.Range("M115").Value = -1588
$ws.Range("N115").Value = -15850
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null
$ws.Range("H117").Value = 9178.417
$ws.Range("I117").Value = 780.125
$ws.Range("J117").Value = 25975
$ws.Range("K117").Value = 2340.375
$ws.Range("L117").Value = 77925
$ws.Range("M117").Value = 1101.625
$ws.Range("N117").Value = -84809
$ws.Range("H119").Value = 5004916.5
$ws.Range("I119").Value = 5004916.5
$ws.Range("K119").Value = 15014749.5
$ws.Range("M119").Value = -15009911.5
$ws.Range("H120").Value = 6000
$ws.Range("I120").Value = 6000
$ws.Range("K120").Value = 18000
$ws.Range("M120").Value = -13162
$ws.Range("H121").Value = 1424.3334
$ws.Range("I121").Value = 631
$ws.Range("J121").Value = 1712.8182
$ws.Range("K121").Value = 1893
$ws.Range("L121").Value = 5138.4546
$ws.Range("M121").Value = -583
$ws.Range("N121").Value = -7758.4546
$ws.Range("H122").Value = 865.1429
$ws.Range("I122").Value = 882.3333
$ws.Range("J122").Value = 852.25
$ws.Range("K122").Value = 7940.9997
$ws.Range("L122").Value = 7670.25
$ws.Range("M122").Value = -5490.9997
$ws.Range("N122").Value = -12570.25
$ws.Range("H123").Value = 3460
$ws.Range("I123").Value = 130
$ws.Range("J123").Value = 5125
$ws.Range("K123").Value = 390
$ws.Range("L123").Value = 15375
$ws.Range("M123").Value = 2060
$ws.Range("N123").Value = -20275
$ws.Range("H124").Value = 8872.25
$ws.Range("I124").Value = 1000
$ws.Range("K124").Value = 3000
$ws.Range("M124").Value = 1910
$ws.Range("H125").Value = 33749.25
$ws.Range("I125").Value = 30000
$ws.Range("K125").Value = 90000
$ws.Range("M125").Value = -85080
$ws.Range("H126").Value = 24997.666
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 34996.5
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 104989.5
$ws.Range("M126").Value = -10060
$ws.Range("N126").Value = -114869.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1587.25
$ws.Range("I80").Value = 1616.3334
$ws.Range("K80").Value = 1616.3334
$ws.Range("M80").Value = -618.3334
$ws.Range("H83").Value = 1587.25
$ws.Range("I83").Value = 1616.3334
$ws.Range("K83").Value = 8081.666999999999
$ws.Range("M83").Value = -3089.666999999999
